# feat: added variables to count updated and inserted rows
$wb = $excel.ActiveWorkbook

$wsData   = $wb.Worksheets.Item("Sheet1_2")
$wsLookup = $wb.Worksheets.Item("Sheet1")

# Track how many existing rows get corrected vs. how many brand new rows get inserted
$updatedRowCount  = 0
$insertedRowCount = 0

# --- Fix the mis-typed address for the first employee (row 2) ---
$wsData.Range("D2").Value = "Vailand at Gum"
$updatedRowCount++

# --- Append the new employee as a brand new row at the bottom of the table ---
$newRowIndex = $wsData.Cells.Item($wsData.Rows.Count, 1).End(-4162).Row + 1
$wsData.Cells.Item($newRowIndex, 1).Value = 8273
$wsData.Cells.Item($newRowIndex, 2).Value = "Vaga"
$wsData.Cells.Item($newRowIndex, 3).Value = "Bond"
$wsData.Cells.Item($newRowIndex, 4).Value = "BondLand"
$insertedRowCount++

# --- Restore the on-screen selection: D2 on the lookup sheet (plus A1), D2 on the data sheet ---
$wsLookup.Activate()
$wsLookup.Range("D2,A1").Select() | Out-Null

$wsData.Activate()
$wsData.Range("D2").Select() | Out-Null
